# 23rd May 2022 update: refresh the ShipmentTrackNum / PackageTrackNum
# tracking numbers in rows 2-22 of Sheet1 with a new batch of values.
#
# Column C = ShipmentTrackNum (always set).
# Column D = PackageTrackNum  (only set on the rows where D already mirrors C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (new tracking number, also update column D?)
$updates = @(
    @{ Row = 2;  Value = "320018606316"; Both = $false }
    @{ Row = 3;  Value = "320018606327"; Both = $false }
    @{ Row = 4;  Value = "320018606350"; Both = $false }
    @{ Row = 5;  Value = "320018606371"; Both = $true  }
    @{ Row = 6;  Value = "320018606419"; Both = $true  }
    @{ Row = 7;  Value = "320018606430"; Both = $true  }
    @{ Row = 8;  Value = "320018606463"; Both = $false }
    @{ Row = 9;  Value = "320018606485"; Both = $false }
    @{ Row = 10; Value = "320018606511"; Both = $false }
    @{ Row = 11; Value = "320018606533"; Both = $false }
    @{ Row = 12; Value = "320018606577"; Both = $false }
    @{ Row = 13; Value = "320018606599"; Both = $true  }
    @{ Row = 14; Value = "320018606625"; Both = $true  }
    @{ Row = 15; Value = "320018606647"; Both = $true  }
    @{ Row = 16; Value = "320018606670"; Both = $true  }
    @{ Row = 17; Value = "320018606691"; Both = $true  }
    @{ Row = 18; Value = "320018606739"; Both = $false }
    @{ Row = 19; Value = "320018606750"; Both = $false }
    @{ Row = 20; Value = "320018606783"; Both = $false }
    @{ Row = 21; Value = "320018606809"; Both = $false }
    @{ Row = 22; Value = "320018606831"; Both = $false }
)

foreach ($u in $updates) {
    # Column C (ShipmentTrackNum). Force text formatting first so the
    # numeric-looking string is kept as text (matches the existing shared
    # string column), then restore the default "Normal" style so no
    # lingering number-format / quote-prefix style is left behind.
    $cellC = $ws.Cells.Item($u.Row, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $u.Value
    $cellC.Style = "Normal"

    if ($u.Both) {
        $cellD = $ws.Cells.Item($u.Row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $u.Value
        $cellD.Style = "Normal"
    }
}
